$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "InfoHub" notes in column G / H for Fermenters & Bright Tanks ---
# These two buildings are now part of the brewery and their "Static?" status is N/A
$ws.Range("G5").Value = "Now part of brewery"
$ws.Range("H5").Value = "N/A"

$ws.Range("G7").Value = "Now part of brewery"
$ws.Range("H7").Value = "N/A"

# --- New note in column I for Bar Stools ---
$ws.Range("I9").Value = "but you sit on them!"

# --- New column J: "Potential freed =" prim counts per building ---
$ws.Range("J1").Value = "Potential freed ="
$ws.Range("J2").Value = 14
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 3
$ws.Range("J8").Value = 2
$ws.Range("J21").Value = 20
$ws.Range("J22").Value = 4

# --- New column K: total potential prims freed (blank cells in J excluded) ---
$ws.Range("K1").Formula = "=SUM(J:J)-COUNT(J:J)"

# --- AutoFit the columns whose contents changed width requirements ---
$ws.Columns.Item(7).AutoFit()  | Out-Null
$ws.Columns.Item(9).AutoFit()  | Out-Null
$ws.Columns.Item(10).AutoFit() | Out-Null

# --- Update the active selection to match where editing finished ---
$ws.Range("J21").Select() | Out-Null
